# Adds the Thursday sub-block (N:V) for 25/3/2019 and a brand-new day block
# (rows 46-52, for 1/4/2019) to the "Pays" sheet, per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Row 41: three new name/cost/subtotal triples (N:P, Q:S, T:V)
# ---------------------------------------------------------------
$ws.Range("N41").Value2 = "chill"
$ws.Range("O41").Value2 = 150
$ws.Range("P41").Formula = "=SUM(O41:O42)"
$ws.Range("P41").HorizontalAlignment = -4108
$ws.Range("P41").VerticalAlignment = -4108

$ws.Range("Q41").Value2 = "ff"
$ws.Range("R41").Value2 = 104
$ws.Range("S41").Formula = "=SUM(R41:R41)"
$ws.Range("S41").HorizontalAlignment = -4108
$ws.Range("S41").VerticalAlignment = -4108

$ws.Range("T41").Value2 = "chill"
$ws.Range("U41").Value2 = 80
$ws.Range("V41").Formula = "=SUM(U41:U41)"
$ws.Range("V41").HorizontalAlignment = -4108
$ws.Range("V41").VerticalAlignment = -4108

# ---------------------------------------------------------------
# Row 42: second line of the new "chill" (N/O) group
# ---------------------------------------------------------------
$ws.Range("N42").Value2 = "dinned"
$ws.Range("O42").Value2 = 55

# ---------------------------------------------------------------
# Row 46 (was a blank row): start of the 1/4/2019 day block
# ---------------------------------------------------------------
# A46 holds a literal date-like label ("1/4/2019"), not a real date
# serial -- pre-format the cell as Text so the string is stored verbatim.
$ws.Range("A46").NumberFormat = "@"
$ws.Range("A46").Value2 = "1/4/2019"
$ws.Range("A46").Style = "Normal"

$ws.Range("B46").Value2 = "dinner"
$ws.Range("C46").Value2 = 55
$ws.Range("D46").Formula = "=SUM(C46:C50)"
$ws.Range("D46").HorizontalAlignment = -4108
$ws.Range("D46").VerticalAlignment = -4108

$ws.Range("E46").Value2 = "dinner"
$ws.Range("F46").Value2 = 215
$ws.Range("G46").Formula = "=SUM(F46:F52)"
$ws.Range("G46").HorizontalAlignment = -4108
$ws.Range("G46").VerticalAlignment = -4108

$ws.Range("H46").Value2 = "snacks"
$ws.Range("I46").Value2 = 30
$ws.Range("J46").Formula = "=SUM(I46:I48)"
$ws.Range("J46").HorizontalAlignment = -4108
$ws.Range("J46").VerticalAlignment = -4108

$ws.Range("K46").Value2 = "dinner"
$ws.Range("L46").Value2 = 180
$ws.Range("M46").Formula = "=SUM(L46:L48)"
$ws.Range("M46").HorizontalAlignment = -4108
$ws.Range("M46").VerticalAlignment = -4108

# ---------------------------------------------------------------
# Row 47
# ---------------------------------------------------------------
$ws.Range("A47").Value2 = "*"
$ws.Range("B47").Value2 = "metro"
$ws.Range("C47").Value2 = 395
$ws.Range("E47").Value2 = "bananas"
$ws.Range("F47").Value2 = 17
$ws.Range("H47").Value2 = "en. drink"
$ws.Range("I47").Value2 = 120
$ws.Range("K47").Value2 = "no"
$ws.Range("L47").Value2 = 110

# ---------------------------------------------------------------
# Row 48
# ---------------------------------------------------------------
$ws.Range("A48").Value2 = "*"
$ws.Range("B48").Value2 = "breakfast"
$ws.Range("C48").Value2 = 110
$ws.Range("E48").Value2 = "water"
$ws.Range("F48").Value2 = 22
$ws.Range("H48").Value2 = "hygiene"
$ws.Range("I48").Value2 = 30
$ws.Range("K48").Value2 = "ff"
$ws.Range("L48").Value2 = 144

# ---------------------------------------------------------------
# Row 49
# ---------------------------------------------------------------
$ws.Range("A49").Value2 = "*"
$ws.Range("B49").Value2 = "hygiene"
$ws.Range("C49").Value2 = 165
$ws.Range("E49").Value2 = "salmon"
$ws.Range("F49").Value2 = 200

# ---------------------------------------------------------------
# Row 50
# ---------------------------------------------------------------
$ws.Range("A50").Value2 = "*"
$ws.Range("B50").Value2 = "bananas"
$ws.Range("C50").Value2 = 70
$ws.Range("E50").Value2 = "bread"
$ws.Range("F50").Value2 = 80

# ---------------------------------------------------------------
# Row 51
# ---------------------------------------------------------------
$ws.Range("A51").Value2 = "*"
$ws.Range("E51").Value2 = "garlic"
$ws.Range("F51").Value2 = 14

# ---------------------------------------------------------------
# Row 52
# ---------------------------------------------------------------
$ws.Range("A52").Value2 = "*"
$ws.Range("E52").Value2 = "blueberries"
$ws.Range("F52").Value2 = 175

# ---------------------------------------------------------------
# Row 53: stays blank, but make it part of the used range so the
# sheet's dimension grows to A1:V53, matching the appended row.
# ---------------------------------------------------------------
$ws.Range("A53").Value2 = "x"
$ws.Range("A53").ClearContents() | Out-Null

# ---------------------------------------------------------------
# New merged ranges (appended after the pre-existing ones)
# ---------------------------------------------------------------
$ws.Range("P41:P42").Merge() | Out-Null
$ws.Range("S41").Merge() | Out-Null
$ws.Range("V41").Merge() | Out-Null
$ws.Range("D46:D50").Merge() | Out-Null
$ws.Range("G46:G52").Merge() | Out-Null
$ws.Range("J46:J48").Merge() | Out-Null
$ws.Range("M46:M48").Merge() | Out-Null
